$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the department column (C) from the full school name to "ENGLISH"
$ws.Range("C2").Value = "ENGLISH"
$ws.Range("C3").Value = "ENGLISH"

# The "10 weeks" duration value was mistakenly entered in the durationMin
# column (F2) instead of the duration column (E2). Move it over.
$ws.Range("E2").Value = $ws.Range("F2").Value2
$ws.Range("F2").Value = $null
